# Update the "two-digit number divided by one-digit number" drill table
# by replacing the division expressions in the 25 populated cells of the table,
# keeping each run's existing character formatting untouched.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "76÷2=" -> "44÷8="
$cellRange = $t.Cell(1, 1).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "44÷8="

# Row 1, Col 2: "85÷8=" -> "28÷8="
$cellRange = $t.Cell(1, 2).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "28÷8="

# Row 1, Col 3: "80÷4=" -> "77÷6="
$cellRange = $t.Cell(1, 3).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "77÷6="

# Row 1, Col 4: "90÷3=" -> "65÷5="
$cellRange = $t.Cell(1, 4).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "65÷5="

# Row 1, Col 5: "73÷3=" -> "83÷9="
$cellRange = $t.Cell(1, 5).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "83÷9="

# Row 5, Col 1: "75÷3=" -> "43÷5="
$cellRange = $t.Cell(5, 1).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "43÷5="

# Row 5, Col 2: "31÷5=" -> "77÷4="
$cellRange = $t.Cell(5, 2).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "77÷4="

# Row 5, Col 3: "88÷5=" -> "30÷4="
$cellRange = $t.Cell(5, 3).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "30÷4="

# Row 5, Col 4: "41÷9=" -> "24÷9="
$cellRange = $t.Cell(5, 4).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "24÷9="

# Row 5, Col 5: "58÷4=" -> "77÷6="
$cellRange = $t.Cell(5, 5).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "77÷6="

# Row 9, Col 1: "70÷6=" -> "93÷7="
$cellRange = $t.Cell(9, 1).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "93÷7="

# Row 9, Col 2: "28÷8=" -> "18÷6="
$cellRange = $t.Cell(9, 2).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "18÷6="

# Row 9, Col 3: "26÷6=" -> "90÷2="
$cellRange = $t.Cell(9, 3).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "90÷2="

# Row 9, Col 4: "68÷6=" -> "74÷4="
$cellRange = $t.Cell(9, 4).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "74÷4="

# Row 9, Col 5: "90÷7=" -> "58÷9="
$cellRange = $t.Cell(9, 5).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "58÷9="

# Row 13, Col 1: "64÷9=" -> "29÷4="
$cellRange = $t.Cell(13, 1).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "29÷4="

# Row 13, Col 2: "84÷2=" -> "23÷4="
$cellRange = $t.Cell(13, 2).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "23÷4="

# Row 13, Col 3: "55÷3=" -> "63÷8="
$cellRange = $t.Cell(13, 3).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "63÷8="

# Row 13, Col 4: "17÷9=" -> "83÷6="
$cellRange = $t.Cell(13, 4).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "83÷6="

# Row 13, Col 5: "74÷5=" -> "31÷5="
$cellRange = $t.Cell(13, 5).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "31÷5="

# Row 17, Col 1: "39÷7=" -> "30÷2="
$cellRange = $t.Cell(17, 1).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "30÷2="

# Row 17, Col 2: "59÷3=" -> "61÷5="
$cellRange = $t.Cell(17, 2).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "61÷5="

# Row 17, Col 3: "95÷8=" -> "45÷6="
$cellRange = $t.Cell(17, 3).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "45÷6="

# Row 17, Col 4: "41÷4=" -> "18÷3="
$cellRange = $t.Cell(17, 4).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "18÷3="

# Row 17, Col 5: "84÷3=" -> "74÷2="
$cellRange = $t.Cell(17, 5).Range
$cellRange.MoveEnd(12, -1) | Out-Null
$cellRange.Text = "74÷2="
